$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New quarter-end columns being appended: R = 31/12/2023, S = 31/03/2024, T = 30/06/2024

# Copy the style from the last existing header cell (Q1: bold/bordered/centered) onto the new header cells
$ws.Range("Q1").Copy($ws.Range("R1:T1"))
$ws.Range("R1").Value = "31/12/2023"
$ws.Range("S1").Value = "31/03/2024"
$ws.Range("T1").Value = "30/06/2024"

# Section-separator rows in the source sheet only contain blank cells across B:Q; replicate that blank
# pattern into R:T by copying an existing blank cell (Q<row>) so the cells exist but stay empty.
$blankRows = @(57, 58, 71, 72, 73, 77, 78)
foreach ($r in $blankRows) {
    $ws.Range("Q$r").Copy($ws.Range("R$r`:T$r"))
}

# Financial data rows: new values for the three added quarters (R, S, T) per row
$newData = @{
    2  = @(798796.032, 757545.9840000001, 720236.992)
    3  = @(454728.992, 418287.008, 378391.008)
    4  = @(152632, 114537, 100844)
    5  = @(0, 0, 0)
    6  = @(102542, 107770, 83105)
    7  = @(64890, 62872, 61457)
    8  = @(0, 0, 0)
    9  = @(97414, 98881, 99461)
    10  = @(0, 0, 0)
    11  = @(37251, 34227, 33524)
    12  = @(113637, 115054, 132209)
    13  = @(0, 0, 0)
    14  = @(0, 0, 0)
    15  = @(0, 0, 0)
    16  = @(3713, 3712, 2815)
    17  = @(0, 0, 0)
    18  = @(0, 0, 0)
    19  = @(35687, 0, 49685)
    20  = @(0, 0, 0)
    21  = @(0, 0, 0)
    22  = @(0, 0, 0)
    23  = @(187744, 181367.008, 166214)
    24  = @(42686, 42838, 43423)
    25  = @(0, 0, 0)
    26  = @(798796.032, 757545.9840000001, 720236.992)
    27  = @(195375.008, 169298, 165531.008)
    28  = @(10510, 11904, 14273)
    29  = @(88236, 74123, 65387)
    30  = @(6903, 6845, 7176)
    31  = @(36042, 57465, 29329)
    32  = @(35, 36, 36)
    33  = @(0, 0, 0)
    34  = @(51539, 17016, 47432)
    35  = @(2110, 1909, 1898)
    36  = @(0, 0, 0)
    37  = @(132394, 138215.008, 120099)
    38  = @(0, 113513, 0)
    39  = @(0, 0, 0)
    40  = @(109168, 0, 94121)
    41  = @(63, 0, 0)
    42  = @(0, 0, 0)
    43  = @(23163, 24702, 25978)
    44  = @(0, 0, 0)
    45  = @(0, 0, 0)
    46  = @(0, 0, 0)
    47  = @(471027.008, 450032.992, 434607.008)
    48  = @(1085844.992, 1085844.992, 1085844.992)
    49  = @(10875, 11154, 10761)
    50  = @(0, 0, 0)
    51  = @(0, 0, 0)
    52  = @(-625692.992, -646966.0159999999, -661998.976)
    53  = @(0, 0, 0)
    54  = @(0, 0, 0)
    55  = @(0, 0, 0)
    56  = @(0, 0, 0)
    59  = @(143240.992, 144927.008, 136900)
    60  = @(-80926.992, -78253, -76291)
    61  = @(62314.008, 66674, 60609)
    62  = @(-59391, -66024, -65218)
    63  = @(-25000, -21712, -22707)
    64  = @(-177, -201, -183)
    65  = @(-596, 117, 9988)
    66  = @(0, 0, -535)
    67  = @(0, 0, 0)
    68  = @(1683, -127, 3014)
    69  = @(10955, 8948, 14534)
    70  = @(-9272, -9075, -11520)
    74  = @(-21167, -21273, -15032)
    75  = @(0, 0, 0)
    76  = @(0, 0, 0)
    79  = @(0, 0, 0)
    80  = @(-21167, -21273, -15032)
}

foreach ($r in $newData.Keys) {
    $vals = $newData[$r]
    $ws.Range("R$r").Value = $vals[0]
    $ws.Range("S$r").Value = $vals[1]
    $ws.Range("T$r").Value = $vals[2]
}